# Update the "Förändrad" (Changed) date column (C) for every data row
# (rows 2 through 236) from 2023-09-10 (serial 45179) to 2023-09-11
# (serial 45180), mirroring the automatic daily refresh of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C236").Value = 45180
